# Rename the default sheet to "Employees"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Employees"

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Value = "First Name"
$ws.Cells.Item(1, 3).Value = "Last Name"
$ws.Cells.Item(1, 4).Value = "Department"
$ws.Cells.Item(1, 5).Value = "Phone"
$ws.Cells.Item(1, 6).Value = "Address"
$ws.Cells.Item(1, 7).Value = "Salary"

# First data row (row 2) - numeric-looking values are forced to text
# (leading apostrophe / quote-prefix) so they are stored as plain text,
# matching the migrated-from-text-file source data.
$ws.Cells.Item(2, 1).Value = "'1"
$ws.Cells.Item(2, 2).Value = "Luke"
$ws.Cells.Item(2, 3).Value = "Phillip"
$ws.Cells.Item(2, 4).Value = "Sales"
$ws.Cells.Item(2, 5).Value = "'1232123"
$ws.Cells.Item(2, 6).Value = "1st Address, Miami"
$ws.Cells.Item(2, 7).Value = "'52000"

# Turn the populated range into a proper Excel Table ("ListObject")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:G2"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table"
